# Adds 'Lucid Trading' as the 13th prop firm across three sheets,
# matching the commit 'Add Lucid Trading as 13th firm - complete
# integration with scoring and UI':
#   - Comparison        -> new row 14 (A14:R14)
#   - Rules Deep Dive   -> new row 16 (A16:I16)
#   - Affiliate Data    -> new row 16 (A16:F16)

$wb = $excel.ActiveWorkbook

# --- Sheet: Comparison (new row 14) ---
$ws = $wb.Worksheets.Item("Comparison")
$ws.Range("A14").Value = 'Lucid Trading'
$ws.Range("B14").Value = '$25K-$150K'
$ws.Range("C14").Value = '$60-$221 one-time (often 40-50% off)'
$ws.Range("D14").Value = '6% ($1.5K-$6K)'
$ws.Range("E14").Value = '$1.2K-$2.7K (soft breach, none on 25K)'
$ws.Range("F14").Value = '$1K-$4.5K EOD trailing'
$ws.Range("G14").Value = 'Trailing (EOD)'
$ws.Range("H14").Value = '90/10'
$ws.Range("I14").Value = 'Daily (5 profitable days per cycle)'
$ws.Range("J14").Value = '1 day to pass'
$ws.Range("K14").Value = 'Unlimited'
$ws.Range("L14").Value = '40% (eval), none on LucidFlex funded'
$ws.Range("M14").Value = 'Rithmic, Tradovate, NinjaTrader, Quantower'
$ws.Range("N14").Value = 'None'
$ws.Range("O14").Value = 4.8
$ws.Range("P14").Value = 2025
$ws.Range("Q14").Value = 'USA'
$ws.Range("R14").Value = '~15 min payouts, one-time fee, no activation fee, LucidFlex no DLL/consistency'

# --- Sheet: Rules Deep Dive (new row 16) ---
$ws = $wb.Worksheets.Item("Rules Deep Dive")
$ws.Range("A16").Value = 'Lucid Trading'
$ws.Range("B16").Value = '1-step'
$ws.Range("C16").Value = 'EOD trailing'
$ws.Range("D16").Value = 'Yes'
$ws.Range("E16").Value = 'Yes (swing trading allowed)'
$ws.Range("F16").Value = 'Yes'
$ws.Range("G16").Value = 'Yes (up to 5 funded accounts)'
$ws.Range("H16").Value = 'LucidScale DLL after consistency'
$ws.Range("I16").Value = '2-10 minis depending on account'

# --- Sheet: Affiliate Data (new row 16) ---
$ws = $wb.Worksheets.Item("Affiliate Data")
$ws.Range("A16").Value = 'Lucid Trading'
$ws.Range("B16").Value = 'Yes'
$ws.Range("C16").Value = '$30-80'
$ws.Range("D16").Value = '40-50% off (SOPF, DGT)'
$ws.Range("E16").Value = 'Monthly'
$ws.Range("F16").Value = 'Fast-growing, strong Discord community'
